$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 171
$ws.Cells.Item(171, 1).Value = 111973697
$ws.Cells.Item(171, 2).Value = 88032
$ws.Cells.Item(171, 4).Value = "VU"
$ws.Cells.Item(171, 5).Value = 6276
$ws.Cells.Item(171, 6).Value = "Goliatmusseron"
$ws.Cells.Item(171, 7).Value = "Tricholoma matsutake"
$ws.Cells.Item(171, 8).Value = "(S.Ito & S.Imai) Singer"
$ws.Cells.Item(171, 17).Value = 438216.5943784415
$ws.Cells.Item(171, 18).Value = 6953090.283452681
$ws.Cells.Item(171, 21).Value = "Härjedalen"
$ws.Cells.Item(171, 23).Value = "Vemdalen"
$ws.Cells.Item(171, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 172
$ws.Cells.Item(172, 1).Value = 111973738
$ws.Cells.Item(172, 2).Value = 90652
$ws.Cells.Item(172, 4).Value = "NT"
$ws.Cells.Item(172, 5).Value = 3100
$ws.Cells.Item(172, 6).Value = "Talltaggsvamp"
$ws.Cells.Item(172, 7).Value = "Bankera fuligineoalba"
$ws.Cells.Item(172, 8).Value = "(Schmidt : Fr.) Pouzar"
$ws.Cells.Item(172, 17).Value = 437818.6501005701
$ws.Cells.Item(172, 18).Value = 6953417.270802823
$ws.Cells.Item(172, 21).Value = "Berg"
$ws.Cells.Item(172, 23).Value = "Åsarne"
$ws.Cells.Item(172, 35).Value = "äldre renbetad lingontallskog med lavfläckar på torr moränmark, på gammal uppgrävd vall"

# Row 173
$ws.Cells.Item(173, 1).Value = 111973718
$ws.Cells.Item(173, 2).Value = 90652
$ws.Cells.Item(173, 4).Value = "NT"
$ws.Cells.Item(173, 5).Value = 3100
$ws.Cells.Item(173, 6).Value = "Talltaggsvamp"
$ws.Cells.Item(173, 7).Value = "Bankera fuligineoalba"
$ws.Cells.Item(173, 8).Value = "(Schmidt : Fr.) Pouzar"
$ws.Cells.Item(173, 17).Value = 437913.625653744
$ws.Cells.Item(173, 18).Value = 6953163.457536075
$ws.Cells.Item(173, 21).Value = "Härjedalen"
$ws.Cells.Item(173, 23).Value = "Vemdalen"
$ws.Cells.Item(173, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 174
$ws.Cells.Item(174, 1).Value = 111973803
$ws.Cells.Item(174, 2).Value = 90652
$ws.Cells.Item(174, 4).Value = "NT"
$ws.Cells.Item(174, 5).Value = 3100
$ws.Cells.Item(174, 6).Value = "Talltaggsvamp"
$ws.Cells.Item(174, 7).Value = "Bankera fuligineoalba"
$ws.Cells.Item(174, 8).Value = "(Schmidt : Fr.) Pouzar"
$ws.Cells.Item(174, 17).Value = 437808.5351863222
$ws.Cells.Item(174, 18).Value = 6953100.247750094
$ws.Cells.Item(174, 21).Value = "Härjedalen"
$ws.Cells.Item(174, 23).Value = "Vemdalen"
$ws.Cells.Item(174, 35).Value = "äldre renbetad lingon- och lavtallskog på torr moränmark"

# Row 175
$ws.Cells.Item(175, 1).Value = 111973722
$ws.Cells.Item(175, 2).Value = 88032
$ws.Cells.Item(175, 4).Value = "VU"
$ws.Cells.Item(175, 5).Value = 6276
$ws.Cells.Item(175, 6).Value = "Goliatmusseron"
$ws.Cells.Item(175, 7).Value = "Tricholoma matsutake"
$ws.Cells.Item(175, 8).Value = "(S.Ito & S.Imai) Singer"
$ws.Cells.Item(175, 17).Value = 437913.6481065798
$ws.Cells.Item(175, 18).Value = 6953091.381175105
$ws.Cells.Item(175, 21).Value = "Härjedalen"
$ws.Cells.Item(175, 23).Value = "Vemdalen"
$ws.Cells.Item(175, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 176
$ws.Cells.Item(176, 1).Value = 111973663
$ws.Cells.Item(176, 2).Value = 90660
$ws.Cells.Item(176, 4).Value = "NT"
$ws.Cells.Item(176, 5).Value = 4362
$ws.Cells.Item(176, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(176, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(176, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(176, 17).Value = 438169.8244046976
$ws.Cells.Item(176, 18).Value = 6953366.599956161
$ws.Cells.Item(176, 21).Value = "Härjedalen"
$ws.Cells.Item(176, 23).Value = "Vemdalen"
$ws.Cells.Item(176, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 182
$ws.Cells.Item(182, 1).Value = 111973665
$ws.Cells.Item(182, 2).Value = 88032
$ws.Cells.Item(182, 4).Value = "VU"
$ws.Cells.Item(182, 5).Value = 6276
$ws.Cells.Item(182, 6).Value = "Goliatmusseron"
$ws.Cells.Item(182, 7).Value = "Tricholoma matsutake"
$ws.Cells.Item(182, 8).Value = "(S.Ito & S.Imai) Singer"
$ws.Cells.Item(182, 17).Value = 438214.6362894689
$ws.Cells.Item(182, 18).Value = 6953402.942781798
$ws.Cells.Item(182, 21).Value = "Härjedalen"
$ws.Cells.Item(182, 23).Value = "Vemdalen"
$ws.Cells.Item(182, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 183
$ws.Cells.Item(183, 1).Value = 111973655
$ws.Cells.Item(183, 2).Value = 90660
$ws.Cells.Item(183, 4).Value = "NT"
$ws.Cells.Item(183, 5).Value = 4362
$ws.Cells.Item(183, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(183, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(183, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(183, 17).Value = 438489.0756873108
$ws.Cells.Item(183, 18).Value = 6953526.341858554
$ws.Cells.Item(183, 21).Value = "Härjedalen"
$ws.Cells.Item(183, 23).Value = "Vemdalen"
$ws.Cells.Item(183, 35).Value = "äldre fattigristallskog på torr moränmark"

# Row 184
$ws.Cells.Item(184, 1).Value = 111973766
$ws.Cells.Item(184, 2).Value = 90652
$ws.Cells.Item(184, 4).Value = "NT"
$ws.Cells.Item(184, 5).Value = 3100
$ws.Cells.Item(184, 6).Value = "Talltaggsvamp"
$ws.Cells.Item(184, 7).Value = "Bankera fuligineoalba"
$ws.Cells.Item(184, 8).Value = "(Schmidt : Fr.) Pouzar"
$ws.Cells.Item(184, 17).Value = 437707.7139296347
$ws.Cells.Item(184, 18).Value = 6953238.955457177
$ws.Cells.Item(184, 21).Value = "Härjedalen"
$ws.Cells.Item(184, 23).Value = "Vemdalen"
$ws.Cells.Item(184, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 185
$ws.Cells.Item(185, 1).Value = 111973711
$ws.Cells.Item(185, 2).Value = 90660
$ws.Cells.Item(185, 4).Value = "NT"
$ws.Cells.Item(185, 5).Value = 4362
$ws.Cells.Item(185, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(185, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(185, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(185, 17).Value = 438030.8777618795
$ws.Cells.Item(185, 18).Value = 6953140.134856743
$ws.Cells.Item(185, 21).Value = "Härjedalen"
$ws.Cells.Item(185, 23).Value = "Vemdalen"
$ws.Cells.Item(185, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 186
$ws.Cells.Item(186, 1).Value = 111973678
$ws.Cells.Item(186, 2).Value = 90654
$ws.Cells.Item(186, 4).Value = "VU"
$ws.Cells.Item(186, 5).Value = 149
$ws.Cells.Item(186, 6).Value = "Tallgråticka"
$ws.Cells.Item(186, 7).Value = "Boletopsis grisea"
$ws.Cells.Item(186, 8).Value = "(Peck) Bondartsev & Singer"
$ws.Cells.Item(186, 17).Value = 438245.3208427017
$ws.Cells.Item(186, 18).Value = 6953249.503443779
$ws.Cells.Item(186, 21).Value = "Härjedalen"
$ws.Cells.Item(186, 23).Value = "Vemdalen"
$ws.Cells.Item(186, 35).Value = "äldre renbetad fattigristallskog med lavfläck på torr moränmark"

# Row 187
$ws.Cells.Item(187, 1).Value = 111973698
$ws.Cells.Item(187, 2).Value = 90660
$ws.Cells.Item(187, 4).Value = "NT"
$ws.Cells.Item(187, 5).Value = 4362
$ws.Cells.Item(187, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(187, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(187, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(187, 17).Value = 438207.123851296
$ws.Cells.Item(187, 18).Value = 6953100.10165237
$ws.Cells.Item(187, 21).Value = "Härjedalen"
$ws.Cells.Item(187, 23).Value = "Vemdalen"
$ws.Cells.Item(187, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 193
$ws.Cells.Item(193, 1).Value = 111973669
$ws.Cells.Item(193, 2).Value = 90660
$ws.Cells.Item(193, 4).Value = "NT"
$ws.Cells.Item(193, 5).Value = 4362
$ws.Cells.Item(193, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(193, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(193, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(193, 17).Value = 437912.8137109271
$ws.Cells.Item(193, 18).Value = 6953242.433193879
$ws.Cells.Item(193, 21).Value = "Härjedalen"
$ws.Cells.Item(193, 23).Value = "Vemdalen"
$ws.Cells.Item(193, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 194
$ws.Cells.Item(194, 1).Value = 111973733
$ws.Cells.Item(194, 2).Value = 90660
$ws.Cells.Item(194, 4).Value = "NT"
$ws.Cells.Item(194, 5).Value = 4362
$ws.Cells.Item(194, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(194, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(194, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(194, 17).Value = 437876.3203048867
$ws.Cells.Item(194, 18).Value = 6953355.130729643
$ws.Cells.Item(194, 21).Value = "Härjedalen"
$ws.Cells.Item(194, 23).Value = "Vemdalen"
$ws.Cells.Item(194, 35).Value = "äldre renbetad lingontallskog med lavfläckar på torr moränmark"

# Row 195
$ws.Cells.Item(195, 1).Value = 111973651
$ws.Cells.Item(195, 2).Value = 90682
$ws.Cells.Item(195, 4).Value = "NT"
$ws.Cells.Item(195, 5).Value = 2059
$ws.Cells.Item(195, 6).Value = "Skrovlig taggsvamp"
$ws.Cells.Item(195, 7).Value = "Hydnellum scabrosum"
$ws.Cells.Item(195, 8).Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Cells.Item(195, 17).Value = 438594.9557070844
$ws.Cells.Item(195, 18).Value = 6953584.041166852
$ws.Cells.Item(195, 21).Value = "Härjedalen"
$ws.Cells.Item(195, 23).Value = "Vemdalen"
$ws.Cells.Item(195, 35).Value = "äldre fattigristallskog på torr moränmark"

# Row 197
$ws.Cells.Item(197, 1).Value = 111973716
$ws.Cells.Item(197, 2).Value = 90660
$ws.Cells.Item(197, 4).Value = "NT"
$ws.Cells.Item(197, 5).Value = 4362
$ws.Cells.Item(197, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(197, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(197, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(197, 17).Value = 437942.142677932
$ws.Cells.Item(197, 18).Value = 6953188.629084867
$ws.Cells.Item(197, 21).Value = "Härjedalen"
$ws.Cells.Item(197, 23).Value = "Vemdalen"
$ws.Cells.Item(197, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 198
$ws.Cells.Item(198, 1).Value = 111973714
$ws.Cells.Item(198, 2).Value = 90658
$ws.Cells.Item(198, 4).Value = "NT"
$ws.Cells.Item(198, 5).Value = 4361
$ws.Cells.Item(198, 6).Value = "Orange taggsvamp"
$ws.Cells.Item(198, 7).Value = "Hydnellum aurantiacum"
$ws.Cells.Item(198, 8).Value = "(Batsch:Fr.) P.Karst."
$ws.Cells.Item(198, 17).Value = 437974.9333715859
$ws.Cells.Item(198, 18).Value = 6953221.065410748
$ws.Cells.Item(198, 21).Value = "Härjedalen"
$ws.Cells.Item(198, 23).Value = "Vemdalen"
$ws.Cells.Item(198, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 199
$ws.Cells.Item(199, 1).Value = 111973744
$ws.Cells.Item(199, 2).Value = 90660
$ws.Cells.Item(199, 4).Value = "NT"
$ws.Cells.Item(199, 5).Value = 4362
$ws.Cells.Item(199, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(199, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(199, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(199, 17).Value = 437686.9991506658
$ws.Cells.Item(199, 18).Value = 6953383.491049132
$ws.Cells.Item(199, 21).Value = "Berg"
$ws.Cells.Item(199, 23).Value = "Åsarne"
$ws.Cells.Item(199, 35).Value = "äldre renbetad lingontallskog med lavfläck på torr moränmark"

# Row 200
$ws.Cells.Item(200, 1).Value = 111973763
$ws.Cells.Item(200, 2).Value = 90652
$ws.Cells.Item(200, 4).Value = "NT"
$ws.Cells.Item(200, 5).Value = 3100
$ws.Cells.Item(200, 6).Value = "Talltaggsvamp"
$ws.Cells.Item(200, 7).Value = "Bankera fuligineoalba"
$ws.Cells.Item(200, 8).Value = "(Schmidt : Fr.) Pouzar"
$ws.Cells.Item(200, 17).Value = 437630.1355663574
$ws.Cells.Item(200, 18).Value = 6953220.681589473
$ws.Cells.Item(200, 21).Value = "Härjedalen"
$ws.Cells.Item(200, 23).Value = "Vemdalen"
$ws.Cells.Item(200, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 201
$ws.Cells.Item(201, 1).Value = 111973764
$ws.Cells.Item(201, 2).Value = 90660
$ws.Cells.Item(201, 4).Value = "NT"
$ws.Cells.Item(201, 5).Value = 4362
$ws.Cells.Item(201, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(201, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(201, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(201, 17).Value = 437776.3383109252
$ws.Cells.Item(201, 18).Value = 6953123.809258236
$ws.Cells.Item(201, 21).Value = "Härjedalen"
$ws.Cells.Item(201, 23).Value = "Vemdalen"
$ws.Cells.Item(201, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 202
$ws.Cells.Item(202, 1).Value = 111973725
$ws.Cells.Item(202, 2).Value = 90660
$ws.Cells.Item(202, 4).Value = "NT"
$ws.Cells.Item(202, 5).Value = 4362
$ws.Cells.Item(202, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(202, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(202, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(202, 17).Value = 437859.825927439
$ws.Cells.Item(202, 18).Value = 6953089.64020811
$ws.Cells.Item(202, 21).Value = "Härjedalen"
$ws.Cells.Item(202, 23).Value = "Vemdalen"
$ws.Cells.Item(202, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 203
$ws.Cells.Item(203, 1).Value = 111973720
$ws.Cells.Item(203, 2).Value = 90652
$ws.Cells.Item(203, 4).Value = "NT"
$ws.Cells.Item(203, 5).Value = 3100
$ws.Cells.Item(203, 6).Value = "Talltaggsvamp"
$ws.Cells.Item(203, 7).Value = "Bankera fuligineoalba"
$ws.Cells.Item(203, 8).Value = "(Schmidt : Fr.) Pouzar"
$ws.Cells.Item(203, 17).Value = 437936.1394948753
$ws.Cells.Item(203, 18).Value = 6953113.912172817
$ws.Cells.Item(203, 21).Value = "Härjedalen"
$ws.Cells.Item(203, 23).Value = "Vemdalen"
$ws.Cells.Item(203, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 204
$ws.Cells.Item(204, 1).Value = 111973689
$ws.Cells.Item(204, 2).Value = 90660
$ws.Cells.Item(204, 4).Value = "NT"
$ws.Cells.Item(204, 5).Value = 4362
$ws.Cells.Item(204, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(204, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(204, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(204, 17).Value = 438478.4845183876
$ws.Cells.Item(204, 18).Value = 6953007.801968225
$ws.Cells.Item(204, 21).Value = "Härjedalen"
$ws.Cells.Item(204, 23).Value = "Vemdalen"
$ws.Cells.Item(204, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"

# Row 205
$ws.Cells.Item(205, 1).Value = 111973658
$ws.Cells.Item(205, 2).Value = 90652
$ws.Cells.Item(205, 4).Value = "NT"
$ws.Cells.Item(205, 5).Value = 3100
$ws.Cells.Item(205, 6).Value = "Talltaggsvamp"
$ws.Cells.Item(205, 7).Value = "Bankera fuligineoalba"
$ws.Cells.Item(205, 8).Value = "(Schmidt : Fr.) Pouzar"
$ws.Cells.Item(205, 17).Value = 438250.4063445947
$ws.Cells.Item(205, 18).Value = 6953324.693784647
$ws.Cells.Item(205, 21).Value = "Härjedalen"
$ws.Cells.Item(205, 23).Value = "Vemdalen"
$ws.Cells.Item(205, 35).Value = "äldre fattigristallskog på torr moränmark"

# Row 206
$ws.Cells.Item(206, 1).Value = 111973660
$ws.Cells.Item(206, 2).Value = 90660
$ws.Cells.Item(206, 4).Value = "NT"
$ws.Cells.Item(206, 5).Value = 4362
$ws.Cells.Item(206, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(206, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(206, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(206, 17).Value = 438117.6502478332
$ws.Cells.Item(206, 18).Value = 6953329.936377568
$ws.Cells.Item(206, 21).Value = "Härjedalen"
$ws.Cells.Item(206, 23).Value = "Vemdalen"
$ws.Cells.Item(206, 35).Value = "äldre fattigristallskog på torr moränmark"

# Row 207
$ws.Cells.Item(207, 1).Value = 111973713
$ws.Cells.Item(207, 2).Value = 90660
$ws.Cells.Item(207, 4).Value = "NT"
$ws.Cells.Item(207, 5).Value = 4362
$ws.Cells.Item(207, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(207, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(207, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(207, 17).Value = 438002.4574124058
$ws.Cells.Item(207, 18).Value = 6953193.462733216
$ws.Cells.Item(207, 21).Value = "Härjedalen"
$ws.Cells.Item(207, 23).Value = "Vemdalen"
$ws.Cells.Item(207, 35).Value = "äldre renbetad fattigris- och lavtallskog på torr moränmark"
